# Project Plan update: reschedule GDD/RTM work, reshuffle CDD tasks for
# Buzzer/DIO, and add new "Review CDD" follow-up tasks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Phase 1: stash a few existing local formats (found on rows 23-28) into
# scratch cells far off to the side, before we overwrite that block, so
# they can still be re-applied afterwards to the handful of cells that
# keep that distinct look.
# ---------------------------------------------------------------------
$ws.Range("A23").Copy()
$ws.Range("ZZ1").PasteSpecial($xlPasteFormats)

$ws.Range("C23").Copy()
$ws.Range("ZZ2").PasteSpecial($xlPasteFormats)

$ws.Range("D23").Copy()
$ws.Range("ZZ3").PasteSpecial($xlPasteFormats)

$ws.Range("H23").Copy()
$ws.Range("ZZ4").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Phase 2: re-base the whole A23:H30 block on row 22's formatting (the
# same look used a few rows earlier in the table).
# ---------------------------------------------------------------------
$ws.Range("A22:H22").Copy()
$ws.Range("A23:H30").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Phase 3: re-apply the handful of cells that keep the older, distinct
# formatting instead of row 22's.
# ---------------------------------------------------------------------
$ws.Range("ZZ1").Copy()
$ws.Range("B26").PasteSpecial($xlPasteFormats)
$ws.Range("F27:F28").PasteSpecial($xlPasteFormats)
$ws.Range("A29:B30").PasteSpecial($xlPasteFormats)

$ws.Range("ZZ2").Copy()
$ws.Range("C27:C28").PasteSpecial($xlPasteFormats)
$ws.Range("C29:C30").PasteSpecial($xlPasteFormats)

$ws.Range("ZZ3").Copy()
$ws.Range("E25:E26").PasteSpecial($xlPasteFormats)
$ws.Range("D29:E30").PasteSpecial($xlPasteFormats)

$ws.Range("ZZ4").Copy()
$ws.Range("H25").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# clean up the scratch cells
$ws.Range("ZZ1:ZZ4").Clear()

# ---------------------------------------------------------------------
# Phase 4: update cell values/content.
# ---------------------------------------------------------------------

# Row 25: "Update GDD requirements" task now finishes a day later and is Done
$ws.Range("E25").Value = 43894
$ws.Range("H25").Value = "Done"

# Row 26: becomes the new "Update RTM with GDD requirements" task
$ws.Range("A26").Value = "Moamen Ahmed"
$ws.Range("B26").Value = "Update RTM with GDD requirements"
$ws.Range("C26").Value = "Non-Technical "
$ws.Range("D26").Value = 43893
$ws.Range("E26").Value = 43894
$ws.Range("F26").Value = "1 day"
$ws.Range("H26").Value = "Pending"

# Row 27: becomes "Create CDD for Buzzer" (now takes 2 days)
$ws.Range("A27").Value = "Areej Helal"
$ws.Range("B27").Value = "Create CDD for Buzzer"
$ws.Range("C27").Value = "Non-Technical "
$ws.Range("D27").Value = 43893
$ws.Range("E27").Value = 43895
$ws.Range("F27").Value = "2 days"
$ws.Range("H27").Value = "Pending"

# Row 28: becomes "Create CDD for DIO" (now takes 2 days)
$ws.Range("A28").Value = "May Abdelsalam"
$ws.Range("B28").Value = "Create CDD for DIO"
$ws.Range("C28").Value = "Non-Technical "
$ws.Range("D28").Value = 43893
$ws.Range("E28").Value = 43895
$ws.Range("F28").Value = "2 days"
$ws.Range("H28").Value = "Pending"

# Row 29: new "Review Buzzer CDD document" task
$ws.Range("A29").Value = "Bishoy Nabil"
$ws.Range("B29").Value = "Review Buzzer CDD document"
$ws.Range("C29").Value = "Non-Technical "
$ws.Range("D29").Value = 43895
$ws.Range("E29").Value = 43895
$ws.Range("F29").Value = "1 day"
$ws.Range("H29").Value = "Pending"

# Row 30: new "Review DIO CDD document" task
$ws.Range("A30").Value = "Mina Yousry"
$ws.Range("B30").Value = "Review DIO CDD document"
$ws.Range("C30").Value = "Non-Technical "
$ws.Range("D30").Value = 43895
$ws.Range("E30").Value = 43895
$ws.Range("F30").Value = "1 day"
$ws.Range("H30").Value = "Pending"
